$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 13, pushing existing rows 13:41 down to 14:42.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new weekly record.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R are constant across the sheet,
# so copy them from the (now shifted) row 14 which held the old row-13 data.
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Vega Monumental Concepción"
$ws.Range("C13").Value = "Bíobío"
$ws.Range("D13").Value = 44453
$ws.Range("E13").Value = 8
$ws.Range("F13").Value = 100112012
$ws.Range("G13").Value = "Espinaca"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 6500
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = 6750
$ws.Range("N13").Value = "$/cuna 10 kilos"
$ws.Range("O13").Value = "Región Metropolitana"
$ws.Range("P13").Value = 675
$ws.Range("Q13").Value = 10
$ws.Range("R13").Value = "Hortaliza"

# Match the date cell style used by the other rows in column D.
$ws.Range("D13").NumberFormat = $ws.Range("D14").NumberFormat
